$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- Capture the "IN PROGRESS" look (fillId 6, centered) from F40 before F40 itself is repainted ---
$ws.Range("F40").Copy()
foreach ($addr in @("F42", "G10", "G11", "G12", "G13", "G14")) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

$ws.Range("F42").Value2 = "IN PROGRESS"
$ws.Range("G10").Value2 = "IN PROGRESS"
$ws.Range("G11").Value2 = "IN PROGRESS"
$ws.Range("G12").Value2 = "IN PROGRESS"
$ws.Range("G13").Value2 = "IN PROGRESS"
$ws.Range("G14").Value2 = "IN PROGRESS"

# --- Copy the "DEFER" look (fillId 5, centered) from F5 onto the new G5/G8 cells ---
$ws.Range("F5").Copy()
foreach ($addr in @("G5", "G8")) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}
$ws.Range("G5").Value2 = "DEFER"
$ws.Range("G8").Value2 = "DEFER"

# --- Row 40 (ProvideMessageArguments): status IN PROGRESS -> DONE (style like other DONE cells, e.g. F2) ---
$ws.Range("F2").Copy()
$ws.Range("F40").PasteSpecial($xlPasteFormats)
$ws.Range("F40").Value2 = "DONE"

# --- Row 40: remove the "Eddy" note entirely ---
$ws.Range("H40").Clear()

# --- Row 42 (ProvideVersionControlProvenance): add the "Eddy" note with plain style (like A42) ---
$ws.Range("A42").Copy()
$ws.Range("H42").PasteSpecial($xlPasteFormats)
$ws.Range("H42").Value2 = "Eddy"

# --- H52 note: drop the rich "TODO: disable the unit test..." run, keep plain "Disabled by default." ---
$ws.Range("H52").Value2 = "Disabled by default."

# --- Selection / scroll state ---
$ws.Range("G14").Select()
